# Update cryptocurrency price/volume data per upstream refresh (GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force the cell to remain a text value (matches source inlineStr cells),
    # even for numeric-looking strings like "149.00" or "0.770", then restore
    # the default "Normal" style so no stray number-format style sticks around.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "67.630.66"
Set-TextValue "E2" "  +7.59%  "
Set-TextValue "D3" "3.588.53"
Set-TextValue "E3" "  +3.20%  "
Set-TextValue "E4" "  -0.11%  "
Set-TextValue "D5" "417.33"
Set-TextValue "E5" "  +0.31%  "
Set-TextValue "D6" "130.51"
Set-TextValue "E6" "  +0.19%  "
Set-TextValue "E7" "  +3.32%  "
Set-TextValue "D8" "3.579.56"
Set-TextValue "E8" "  +3.16%  "
Set-TextValue "E9" "  -0.09%  "
Set-TextValue "D10" "0.770"
Set-TextValue "E10" "  +5.65%  "
Set-TextValue "D11" "0.178"
Set-TextValue "E11" "  +16.53%  "
Set-TextValue "D12" "0.0000342"
Set-TextValue "E12" "  +51.00%  "
Set-TextValue "D13" "42.41"
Set-TextValue "E13" "  -0.28%  "
Set-TextValue "D14" "10.00"
Set-TextValue "E14" "  +1.88%  "
Set-TextValue "D15" "4.138.47"
Set-TextValue "E15" "  +2.57%  "
Set-TextValue "E16" "  -0.24%  "
Set-TextValue "D17" "20.46"
Set-TextValue "E17" "  -0.65%  "
Set-TextValue "D18" "3.591.85"
Set-TextValue "E18" "  +2.80%  "
Set-TextValue "D19" "1.15"
Set-TextValue "E19" "  +6.14%  "
Set-TextValue "D20" "67.371.68"
Set-TextValue "E20" "  +7.26%  "
Set-TextValue "E21" "  -2.46%  "
Set-TextValue "D22" "458.72"
Set-TextValue "E22" "  -1.11%  "
Set-TextValue "D23" "88.41"
Set-TextValue "E23" "  -2.27%  "
Set-TextValue "B24" "ImmutableX"
Set-TextValue "C24" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D24" "3.10"
Set-TextValue "E24" "  -6.06%  "
Set-TextValue "B25" "InternetComputer(DFINITY)"
Set-TextValue "C25" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D25" "13.47"
Set-TextValue "E25" "  +1.41%  "
Set-TextValue "D26" "3.37"
Set-TextValue "E26" "  +1.55%  "
Set-TextValue "D27" "10.12"
Set-TextValue "E27" "  -6.34%  "
Set-TextValue "D28" "35.04"
Set-TextValue "E28" "  +4.64%  "
Set-TextValue "E29" "  +1.27%  "
Set-TextValue "D30" "2.80"
Set-TextValue "E30" "  +4.48%  "
Set-TextValue "D31" "12.40"
Set-TextValue "E31" "  +1.85%  "
Set-TextValue "B32" "RenderToken"
Set-TextValue "C32" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D32" "7.47"
Set-TextValue "E32" "  -1.38%  "
Set-TextValue "B33" "Hedera"
Set-TextValue "C33" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D33" "0.117"
Set-TextValue "E33" "  +4.19%  "
Set-TextValue "D34" "41.78"
Set-TextValue "E34" "  +1.54%  "
Set-TextValue "D35" "0.163"
Set-TextValue "E35" "  -4.38%  "
Set-TextValue "D36" "0.999"
Set-TextValue "E36" "  -0.03%  "
Set-TextValue "D37" "56.62"
Set-TextValue "E37" "  -2.59%  "
Set-TextValue "E38" "  +0.75%  "
Set-TextValue "E39" "  +21.85%  "
Set-TextValue "E40" "  +8.97%  "
Set-TextValue "E41" "  -0.19%  "
Set-TextValue "E42" "  +0.00%  "
Set-TextValue "D43" "149.00"
Set-TextValue "E43" "  +0.67%  "
Set-TextValue "E44" "  -0.26%  "
Set-TextValue "E45" "  -1.84%  "
Set-TextValue "E46" "  -2.74%  "
Set-TextValue "E47" "  -3.66%  "
Set-TextValue "E48" "  -3.93%  "
Set-TextValue "E49" "  -2.03%  "
Set-TextValue "E50" "  +14.89%  "
Set-TextValue "D51" "15.67"
Set-TextValue "E51" "  -4.55%  "
